$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "56.617.77"
$ws.Range("E2").Value = "  +0.48%  "

# Row 3
$ws.Range("D3").Value = "2.390.11"
$ws.Range("E3").Value = "  +0.62%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.40"
$ws.Range("E5").Value = "  +1.61%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.80"
$ws.Range("E6").Value = "  +4.02%  "

# Row 7
$ws.Range("E7").Value = "  +0.17%  "

# Row 9
$ws.Range("D9").Value = "2.393.09"
$ws.Range("E9").Value = "  -0.11%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0976"
$ws.Range("E10").Value = "  +2.28%  "

# Row 11
$ws.Range("E11").Value = "  +0.36%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.337"
$ws.Range("E12").Value = "  +6.39%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.69"
$ws.Range("E13").Value = "  +0.64%  "

# Row 14
$ws.Range("D14").Value = "2.815.13"
$ws.Range("E14").Value = "  +0.70%  "

# Row 15
$ws.Range("D15").Value = "56.576.82"
$ws.Range("E15").Value = "  +0.66%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.79"
$ws.Range("E16").Value = "  +1.96%  "

# Row 17
$ws.Range("E17").Value = "  +1.57%  "

# Row 18
$ws.Range("D18").Value = "2.359.23"
$ws.Range("E18").Value = "  +1.53%  "

# Row 19
$ws.Range("E19").Value = "  +0.87%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.05"
$ws.Range("E20").Value = "  +0.87%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "310.22"
$ws.Range("E21").Value = "  +0.56%  "

# Row 22
$ws.Range("E22").Value = "  +0.55%  "

# Row 23
$ws.Range("E23").Value = "  +0.51%  "

# Row 24
$ws.Range("E24").Value = "  -1.39%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.45"
$ws.Range("E25").Value = "  +0.94%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.34%  "

# Row 27
$ws.Range("B27").Value = "Polygon"
$ws.Range("C27").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.374"
$ws.Range("E27").Value = "  +0.12%  "

# Row 28
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.150"
$ws.Range("E28").Value = "  +0.40%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.34"
$ws.Range("E29").Value = "  +1.76%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "173.88"
$ws.Range("E30").Value = "  +0.91%  "

# Row 31
$ws.Range("D31").Value = "0.0₃0727"
$ws.Range("E31").Value = "  +2.54%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.66"
$ws.Range("E32").Value = "  +0.42%  "

# Row 33
$ws.Range("E33").Value = "  +1.94%  "

# Row 34
$ws.Range("E34").Value = "  -3.87%  "

# Row 35
$ws.Range("E35").Value = "  +0.13%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("E36").Value = "  +0.30%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.89"
$ws.Range("E37").Value = "  +0.86%  "

# Row 38
$ws.Range("E38").Value = "  -0.12%  "

# Row 39
$ws.Range("E39").Value = "  +1.94%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.58"
$ws.Range("E40").Value = "  +2.02%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.814"
$ws.Range("E41").Value = "  +4.06%  "

# Row 42
$ws.Range("E42").Value = "  +1.11%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "132.20"
$ws.Range("E43").Value = "  +2.69%  "

# Row 44
$ws.Range("E44").Value = "  +2.24%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.85"
$ws.Range("E45").Value = "  +1.91%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.566"
$ws.Range("E46").Value = "  +0.58%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0909"
$ws.Range("E47").Value = "  +1.48%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "247.49"
$ws.Range("E48").Value = "  -1.60%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0485"
$ws.Range("E49").Value = "  +0.63%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0210"
$ws.Range("E50").Value = "  +1.84%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.18"
$ws.Range("E51").Value = "  +6.99%  "
